$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Formula = '''96.793.50'
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.Formula = '''  +0.45%  '
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.Formula = '''3.691.11'
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.Formula = '''  +3.37%  '
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.Formula = '''  -0.05%  '
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.Formula = '''243.90'
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.Formula = '''  +1.39%  '
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.Formula = '''1.89'
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.Formula = '''  +15.90%  '
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.Formula = '''669.40'
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.Formula = '''  +2.50%  '
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.Formula = '''0.430'
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.Formula = '''  +5.19%  '
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.Formula = '''1.11'
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.Formula = '''  +4.88%  '
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.Formula = '''3.687.60'
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.Formula = '''  +3.35%  '
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.Formula = '''45.59'
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.Formula = '''  +5.19%  '
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.Formula = '''  +1.62%  '
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.Formula = '''6.63'
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.Formula = '''  +3.74%  '
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.Formula = '''4.375.59'
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.Formula = '''  +3.29%  '
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.Formula = '''  +4.30%  '
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.Formula = '''96.507.76'
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.Formula = '''  +0.30%  '
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.Formula = '''9.03'
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.Formula = '''  +16.21%  '
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.Formula = '''3.687.46'
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.Formula = '''  +3.29%  '
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.Formula = '''  +1.15%  '
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.Formula = '''18.50'
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.Formula = '''  +3.91%  '
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.Formula = '''0.536'
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.Formula = '''  +3.80%  '
$r.Style = "Normal"

$r = $ws.Range("B23")
$r.Formula = '''SuiNetwork'
$r.Style = "Normal"

$r = $ws.Range("C23")
$r.Formula = '''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.Formula = '''3.51'
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.Formula = '''  +3.39%  '
$r.Style = "Normal"

$r = $ws.Range("B24")
$r.Formula = '''BitcoinCash'
$r.Style = "Normal"

$r = $ws.Range("C24")
$r.Formula = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.Formula = '''522.56'
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.Formula = '''  +3.49%  '
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.Formula = '''0.0000208'
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.Formula = '''  +4.56%  '
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.Formula = '''7.04'
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.Formula = '''  +1.01%  '
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.Formula = '''102.93'
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.Formula = '''  +7.02%  '
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.Formula = '''  +0.03%  '
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.Formula = '''  +7.52%  '
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.Formula = '''  +1.63%  '
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.Formula = '''  +6.89%  '
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.Formula = '''0.998'
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.Formula = '''  -0.23%  '
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.Formula = '''  +1.67%  '
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.Formula = '''  +11.21%  '
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.Formula = '''32.97'
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.Formula = '''  +5.19%  '
$r.Style = "Normal"

$r = $ws.Range("D36")
$r.Formula = '''0.999'
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.Formula = '''  -0.44%  '
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.Formula = '''0.589'
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.Formula = '''  +4.18%  '
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.Formula = '''622.77'
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.Formula = '''8.83'
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.Formula = '''  -0.46%  '
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.Formula = '''42.83'
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.Formula = '''  +29.35%  '
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.Formula = '''0.161'
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.Formula = '''  +6.95%  '
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.Formula = '''0.961'
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.Formula = '''  +6.31%  '
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.Formula = '''  +7.91%  '
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.Formula = '''  +8.54%  '
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.Formula = '''0.0459'
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.Formula = '''  +7.91%  '
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.Formula = '''0.434'
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.Formula = '''  +25.08%  '
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.Formula = '''  +0.92%  '
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.Formula = '''23.63'
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.Formula = '''  +0.44%  '
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.Formula = '''  +4.63%  '
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.Formula = '''54.66'
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.Formula = '''  +3.63%  '
$r.Style = "Normal"
